$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Data ("base de datos") update: "Periodo Mora" column (E16:E21) ---
# Previous periods were listed oldest-to-newest (2412, 2501, 2502, 2503, 2504, 2505).
# Old statements were removed and the list was refreshed; the remaining periods are
# now presented in reverse order for the first six rows while the last row (2506)
# is unchanged.
$ws.Range("E16").Value = "2505"
$ws.Range("E17").Value = "2504"
$ws.Range("E18").Value = "2503"
$ws.Range("E19").Value = "2502"
$ws.Range("E20").Value = "2501"
$ws.Range("E21").Value = "2412"
$ws.Range("E22").Value = "2506"

# --- Column width refresh (best-fit recalculated by a newer Excel build) ---
$ws.Columns.Item(2).ColumnWidth = 17.709635416666668
$ws.Columns.Item(3).ColumnWidth = 15.893229166666666
$ws.Columns.Item(5).ColumnWidth = 12.709635416666666
$ws.Columns.Item(6).ColumnWidth = 9.346354166666666
$ws.Columns.Item(7).ColumnWidth = 13.529947916666666
$ws.Columns.Item(8).ColumnWidth = 18.529947916666668
$ws.Columns.Item(9).ColumnWidth = 17.256510416666668
$ws.Columns.Item(10).ColumnWidth = 14.166666666666666
